# Slide 4 ("2. 팀 구성 및 역할 (Team Project)") - "TextBox 4" shape:
#   1. "□ 팀장 " -> "□ 백종우 팀장 " (team-lead bullet gets the member's name)
#   2. In 안재원's bullet, the " " run and the "발표 " run (both plain
#      ko-KR runs, no bold) get merged into a single " 발표 " run.

$p   = $ppt.ActivePresentation
$s   = $p.Slides.Item(4)
$shp = $s.Shapes.Item(3)
$tr  = $shp.TextFrame.TextRange

# the shape auto-fits its height to the text (<a:spAutoFit/>); remember the
# stored height so it can be restored after the text edit below (this host
# re-lays-out / re-measures spAutoFit shapes whenever their text changes,
# same as PowerPoint does on a live edit - but this particular file's
# original height predates that recompute, so put it back afterwards)
$origHeightPt = $shp.Height

# --- 1) "□ 팀장 " -> "□ 백종우 팀장 " -------------------------------------
$full = $tr.Text
$old1 = "□ 팀장 "
$pos1 = $full.IndexOf($old1)
if ($pos1 -ge 0) {
    $tr.Characters($pos1 + 1, $old1.Length).Text = "□ 백종우 팀장 "
}

# --- 2) merge " " + "발표 " runs on 안재원's line -------------------------
$full = $tr.Text
$marker = "test_cart, test_checkout 제작,"
$pos2 = $full.IndexOf($marker)
if ($pos2 -ge 0) {
    $start2 = $pos2 + $marker.Length + 1
    $tr.Characters($start2, 4).Text = " 발표 "
}

# --- restore the shape's original (pre-edit) height ------------------------
$targetEmu = [math]::Round($origHeightPt * 12700.0)
$shp.Height = $origHeightPt
$bestPt = $origHeightPt
$bestDiff = [math]::Abs([math]::Round($shp.Height * 12700.0) - $targetEmu)
if ($bestDiff -ne 0) {
    for ($i = -500; $i -le 500; $i++) {
        $cand = $origHeightPt + ($i * 0.00001)
        $shp.Height = $cand
        $diff = [math]::Abs([math]::Round($shp.Height * 12700.0) - $targetEmu)
        if ($diff -lt $bestDiff) {
            $bestDiff = $diff
            $bestPt = $cand
        }
        if ($diff -eq 0) {
            break
        }
    }
    $shp.Height = $bestPt
}
